# eventbuttons.xlsx: add four new "Commands" entries (showCurve, showExtraCurve,
# showEvents, showBackgroundEvents) just before the existing "RC Command" block
# on the Commands sheet, shifting everything below down by four rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Commands")

# The sheet carries a trailing placeholder row at the very bottom of the grid
# (row 1048576) with no cell content, only row formatting. Drop it first,
# while it is still addressable, so it does not get dragged along (and pushed
# out of range) by the upcoming row insert.
$ws.Rows.Item(1048576).Delete($null)

# Insert 4 blank rows above the current row 100 (the "RC Command" row),
# pushing that row and everything after it down to rows 104+.
$ws.Range("A100:A103").EntireRow.Insert()

# Populate the newly inserted rows with the new command documentation.
$ws.Range("B100").Value = "showCurve(<name>,<bool>)"
$ws.Range("C100").Value = "shows/hides the curve indicated by <name> which is one of { ET, BT, DeltaET, DeltaBT, BackgroundET, BackgroundBT}"

$ws.Range("B101").Value = "showExtraCurve(<extra_device>,<curve>,<bool>)"
$ws.Range("C101").Value = "shows/hides the <curve> (one of {T1,T2}) of the zero-based <extra_device> number"

$ws.Range("B102").Value = "showEvents(<event_type>, <bool>)"
$ws.Range("C102").Value = "shows/hides the events of <event_type> in [1,..,5]"

$ws.Range("B103").Value = "showBackgroundEvents(<bool>)"
$ws.Range("C103").Value = "shows/hides the events of the background profile"

# Match the author's resulting selection on the Commands sheet.
$ws.Activate() | Out-Null
$ws.Range("C101").Select() | Out-Null
